$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "BUFF_LOAD"
$ws.Range("F11").Value = "Tells the micro to start looking for packets"
$ws.Range("E12").Value = "SHOW_MODE"
$ws.Range("F12").Value = "Tells the mirco to just light up drums that were hit"

$ws.Range("F12").Select()
$excel.ActiveWindow.ScrollColumn = 3
